# Apply the changes described by the commit:
#  1. Refresh the cached "datetime1" field text (date placeholder default
#     text) from 2/25/2020 to 3/3/2020 across every slide layout and the
#     slide master.
#  2. Nudge the picture on slide 6 slightly down.
#  3. Add a new textbox with a hyperlink under the picture on slide 6.

$p = $ppt.ActivePresentation

# --- 1. Update the cached date placeholder text (layouts + master) -------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$layouts = $master.CustomLayouts

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "3/3/2020"
        }
    }
}

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "3/3/2020"
    }
}

# --- 2. Move the picture on slide 6 down slightly -------------------------
$slide6 = $p.Slides.Item(6)
$pic = $slide6.Shapes.Item(2)
$pic.Top = 2655248 / 12700.0

# --- 3. Add the new hyperlink textbox on slide 6 --------------------------
$textBox = $slide6.Shapes.AddTextbox(1, 3476531 / 12700.0, 5846074 / 12700.0, 5238935 / 12700.0, 369332 / 12700.0)
$textBox.TextFrame.TextRange.Text = "http://openaccess.thecvf.com/CVPR2019.py"
$textBox.TextFrame.TextRange.ActionSettings(1).Hyperlink.Address = "http://openaccess.thecvf.com/CVPR2019.py"
$textBox.TextFrame.WordWrap = 0
$textBox.TextFrame.AutoSize = 1
$textBox.Fill.Visible = 0

Write-Host "Edit complete"
